$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = '51.690.63'
$ws.Range("E2").Value = '  -0.48%  '

# Row 3: update D3, E3
$ws.Range("D3").Value = '2.802.76'
$ws.Range("E3").Value = '  +0.38%  '

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '354.71'
$ws.Range("E5").Value = '  -1.15%  '

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.09'
$ws.Range("E6").Value = '  -0.58%  '

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.556'
$ws.Range("E7").Value = '  -0.93%  '

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.02%  '

# Row 9: update E9
$ws.Range("E9").Value = '  +5.58%  '

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.92'
$ws.Range("E10").Value = '  -0.80%  '

# Row 11: update E11
$ws.Range("E11").Value = '  +0.63%  '

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0837'
$ws.Range("E12").Value = '  -1.74%  '

# Row 13: update D13, E13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.99'
$ws.Range("E13").Value = '  +2.18%  '

# Row 14: update E14
$ws.Range("E14").Value = '  +1.45%  '

# Row 15: update D15, E15
$ws.Range("D15").Value = '3.233.33'
$ws.Range("E15").Value = '  +0.15%  '

# Row 16: update D16, E16
$ws.Range("D16").Value = '2.792.79'
$ws.Range("E16").Value = '  -0.09%  '

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.946'
$ws.Range("E17").Value = '  +0.05%  '

# Row 18: update D18, E18
$ws.Range("D18").Value = '51.612.89'
$ws.Range("E18").Value = '  -0.52%  '

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("E19").Value = '  +3.68%  '

# Row 20: update E20
$ws.Range("E20").Value = '  +3.14%  '

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  +1.60%  '

# Row 22: update D22, E22
$ws.Range("D22").Value = '0.0₃0977'
$ws.Range("E22").Value = '  -0.10%  '

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.36'
$ws.Range("E23").Value = '  +0.27%  '

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.98'
$ws.Range("E24").Value = '  -0.84%  '

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.75'
$ws.Range("E25").Value = '  -0.92%  '

# Row 26: update E26
$ws.Range("E26").Value = '  +0.03%  '

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.07'
$ws.Range("E27").Value = '  -1.68%  '

# Row 28: update E28
$ws.Range("E28").Value = '  -0.13%  '

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.34'
$ws.Range("E29").Value = '  +0.78%  '

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.12'
$ws.Range("E30").Value = '  +6.75%  '

# Row 31: update B31, C31, D31, E31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("E31").Value = '  +1.93%  '

# Row 32: update B32, C32, D32, E32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +0.50%  '

# Row 33: update B33, C33, D33, E33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.94'
$ws.Range("E33").Value = '  -0.48%  '

# Row 34: update B34, C34, D34, E34
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  +10.14%  '

# Row 35: update B35, C35, D35, E35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0444'
$ws.Range("E35").Value = '  -4.65%  '

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0856'
$ws.Range("E36").Value = '  +0.79%  '

# Row 37: update E37
$ws.Range("E37").Value = '  -0.12%  '

# Row 38: update D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.91'
$ws.Range("E38").Value = '  +0.46%  '

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.00'
$ws.Range("E39").Value = '  +0.86%  '

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.13'
$ws.Range("E40").Value = '  -2.58%  '

# Row 41: update E41
$ws.Range("E41").Value = '  +0.62%  '

# Row 42: update D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.48'
$ws.Range("E42").Value = '  -4.17%  '

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.54'
$ws.Range("E43").Value = '  +0.17%  '

# Row 44: update E44
$ws.Range("E44").Value = '  -1.53%  '

# Row 45: update D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.77'
$ws.Range("E45").Value = '  -1.19%  '

# Row 46: update D46, E46
$ws.Range("D46").Value = '2.129.52'
$ws.Range("E46").Value = '  +2.26%  '

# Row 47: update B47, C47, D47, E47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  +6.81%  '

# Row 48: update B48, C48, D48, E48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").Value = '  +2.21%  '

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.911'
$ws.Range("E49").Value = '  -4.51%  '

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.39'
$ws.Range("E50").Value = '  -5.82%  '

# Row 51: update E51
$ws.Range("E51").Value = '  +7.07%  '

